$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "pressen","werfen","locken","treiben","altern","achten","sammeln","landen","wirken","zeigen",
    "binden","kosten","mögen","treten","jubeln","schauen","schicken","bergen","lesen","ziehen",
    "schwächen","sterben","pissen","fischen","spinnen","decken","sprengen","gelten","nutzen","boxen",
    "schwingen","ärgern","hindern","reiten","albern","stärken","heilen","dürfen","schreiten","kichern",
    "münzen","löschen","stehlen","helfen","scheiden","melden","zünden","stechen","orten","seufzen",
    "flüchten","plaudern","flehen","loben","feuern","filmen","zeichnen","betteln","fällen","heben",
    "fließen","klingen","bellen","bilden","wehtun","saufen","duschen","pflanzen","arten","freuen",
    "hoffen","tropfen","dauern","scheitern","räumen","rauchen","sinken","zielen","wehen","bluten",
    "wenden","dringen","herrschen","knurren","stoßen","wüten","heulen","lügen","jagen","planen",
    "formen","deuten","drehen","platzen","mauern","beißen","trennen","sagen","drohen","sichern",
    "schmecken","spüren","reizen","feiern","irren","schlafen","führen","rasen","rechnen","trauen",
    "liefern","fühlen","zählen","malen","klagen","scheinen","sorgen","weichen","herzen","ändern",
    "schulden","holen","gründen","geben","machen","schrecken","wundern","flüstern","zögern","dienen",
    "warnen","runden","quälen","tauschen","grüßen","tollen","fallen","hören","kehren","starren",
    "biegen","erben","siegen","kriegen","kämpfen","graben","kümmern","suchen","schenken","brauchen",
    "wärmen","pfeifen","bitten","spielen","krachen","greifen","enden","hauen","rufen","backen",
    "wechseln","lohnen","streichen","öffnen","äußern","bieten","schultern","laufen","stecken","sperren",
    "passen","regeln","töten","werden","süßen","folgen","schlucken","ehren","zahlen","klettern",
    "kürzen","liegen","stammen","lockern","wachsen","fahren","lieben","knarren","schwören","fangen",
    "bauen","stürmen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
